$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the footnote marker symbols (superscript asterisk/dagger/etc.) to
# superscript letters in the row labels.
$ws.Range("A4").Value = "        Organophosphateᵃ"
$ws.Range("A5").Value = "        Pyrethroidᵇ"
$ws.Range("A7").Value = "        Carbamateᶜ"
$ws.Range("A15").Value = "        Triazineᵈ"
$ws.Range("A60").Value = "    Wood Preservativesᵉ"

# Update the footnote legend text block: replace the symbol markers with
# superscript letters, and drop the blank line that used to separate the
# footnote list from the "Abbreviations" line.
$footnote = "ᵃ Includes the organophosphate breakdown product, diethyl phosphate" + [char]10 + `
            "ᵇ Includes the pyrethroid breakdown product, 3-phenoxybenzoic acid" + [char]10 + `
            "ᶜ Includes the carbamate breakdown product, aldicarb sulfone" + [char]10 + `
            "ᵈ Includes a triazinone (metribuzin)" + [char]10 + `
            "ᵉ Includes the wood preservative breakdown product, pentachloroanisole" + [char]10 + `
            "Abbreviations: UV = Ultraviolet"
$ws.Range("A66").Value = $footnote
